$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I and J, rows 2-18 ([row] = I, J)
$data = @{
    2  = @(1, 5)
    3  = @(1, 3)
    4  = @(1, 5)
    5  = @(1, 5)
    6  = @(1, 2)
    7  = @(5, 8)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 4)
    11 = @(1, 5)
    12 = @(1, 4)
    13 = @(1, 4)
    14 = @(1, 5)
    15 = @(5, 9)
    16 = @(1, 4)
    17 = @(1, 3)
    18 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
